# DailyWorkReport.xlsx - add the "Versioning and Basic web development final
# demo" day blocks (rows 67-74) following the same layout used by the
# preceding weeks in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 67: new day header (date + "Domm") -------------------------------
$ws.Range("A54").Copy()
$ws.Range("A67").PasteSpecial(-4122)
$ws.Range("B54:D54").Copy()
$ws.Range("B67").PasteSpecial(-4122)

$ws.Range("A67").Value = 45680
$ws.Range("B67").Value = "Domm"
$ws.Range("D67").Value = 0.25

# --- row 68: Meeting / Reconsile -------------------------------------------
$ws.Range("A55:D55").Copy()
$ws.Range("A68").PasteSpecial(-4122)

$ws.Range("B68").Value = "Meeting"
$ws.Range("C68").Value = "Reconsile"
$ws.Range("D68").Value = 1

# --- row 69: General Discussion(offline + online) --------------------------
$ws.Range("A57:D57").Copy()
$ws.Range("A69").PasteSpecial(-4122)

$ws.Range("C69").Value = "General Discussion(offline + online)"
$ws.Range("D69").Value = 0.75

# --- row 70: Study / Web Development Revision for Reconsile ----------------
$ws.Range("A58:D58").Copy()
$ws.Range("A70").PasteSpecial(-4122)

$ws.Range("B70").Value = "Study"
$ws.Range("C70").Value = "Web Development Revision for Reconsile"
$ws.Range("D70").Value = 2

# --- row 71: Debugging Revision for Reconsile -------------------------------
$ws.Range("A59:D59").Copy()
$ws.Range("A71").PasteSpecial(-4122)

$ws.Range("C71").Value = "Debugging Revision for Reconsile"
$ws.Range("D71").Value = 0.5

# --- row 72: .Net Core Fundamental ------------------------------------------
$ws.Range("A60:D60").Copy()
$ws.Range("A72").PasteSpecial(-4122)

$ws.Range("C72").Value = ".Net Core Fundamental"
$ws.Range("D72").Value = 1

# --- row 73: Middleware ------------------------------------------------------
$ws.Range("A60:D60").Copy()
$ws.Range("A73").PasteSpecial(-4122)

$ws.Range("C73").Value = "Middleware"
$ws.Range("D73").Value = 2.5

# --- row 74: trailing blank row, borders cleared ----------------------------
$ws.Range("A65:D65").Copy()
$ws.Range("A74").PasteSpecial(-4122)
$ws.Range("A74:D74").Borders.LineStyle = -4142

$excel.CutCopyMode = $false

# --- sheet view: scroll to the new block, select the new block -------------
$ws.Range("A67:D74").Select()
$excel.ActiveWindow.ScrollRow = 46

